$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Helper pattern used throughout:
#   - To change the *content* of an existing paragraph while keeping it in
#     place, we replace that paragraph's whole Range with a freshly built
#     <w:p> (this lets us control run-splitting / proofErr wrappers exactly).
#   - To *add* new paragraphs after an existing one, we first call
#     InsertParagraphAfter() to create an empty placeholder paragraph, then
#     InsertXML() on that placeholder's Range with one-or-more <w:p> nodes.
# We work from the bottom of the document upward so that paragraph indices
# above the paragraph currently being edited never shift under us.
# ---------------------------------------------------------------------------

# --- A) Paragraph 23: "Conclusions and Future work:" ----------------------
$p23 = $d.Paragraphs(23)
$xmlA = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">Expanded </w:t></w:r><w:r><w:t>Stock assessment</w:t></w:r><w:r><w:t xml:space="preserve"> paragraph</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$p23.Range.InsertXML($xmlA)

$p23b = $d.Paragraphs(23)
$p23b.Range.InsertParagraphAfter()
$newConclusions = $d.Paragraphs(24)
$xmlA2 = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Conclusions and Future work:</w:t></w:r><w:r><w:t xml:space="preserve"> Broad EBFM paragraph</w:t></w:r></w:p>'
$newConclusions.Range.InsertXML($xmlA2)

# --- B) Paragraph 22: "Limitations of our study:" (unchanged) -------------
#     add a new ilvl-1 paragraph right after it
$p22 = $d.Paragraphs(22)
$p22.Range.InsertParagraphAfter()
$newOnlyAssessment = $d.Paragraphs(23)
$xmlB = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Only assessment models </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>reports</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> (not control rules or reference points)</w:t></w:r></w:p>'
$newOnlyAssessment.Range.InsertXML($xmlB)

# --- C) Paragraph 21: "How to better track implicit influences? Council docs"
$p21 = $d.Paragraphs(21)
$xmlC = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">What are other ways EBFM can advance (use ecosystem </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>info)</w:t></w:r><w:r><w:t>How</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> to better track implicit influences? Council docs</w:t></w:r></w:p>'
$p21.Range.InsertXML($xmlC)

# --- D) Paragraph 20: "Already sampling effects directly (compositional data)"
$p20 = $d.Paragraphs(20)
$xmlD = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Already sampling effects directly (compositional data)</w:t></w:r></w:p>'
$p20.Range.InsertXML($xmlD)

# --- E) Paragraph 19: "When is more information not helpful?" -> unchanged

# --- F) Paragraph 18: "What are the barriers?" -----------------------------
$p18 = $d.Paragraphs(18)
$xmlF = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>What are the barriers</w:t></w:r><w:r><w:t xml:space="preserve"> to more progress</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p>'
$p18.Range.InsertXML($xmlF)

$p18b = $d.Paragraphs(18)
$p18b.Range.InsertParagraphAfter()
$newBlock1 = $d.Paragraphs(19)
$xmlF2 = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Overfished status</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Climate and fishing combined</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Data</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Funding for diet studies</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>expertise</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Institutional (stock assessment process)</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Barriers</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Creative process</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Solutions</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Ecosystem and assessment working together</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Moving assessment authors </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>around ?</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> (exchange)</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr></w:p>'
$newBlock1.Range.InsertXML($xmlF2)

# --- G) Paragraph 17: "Summary of results- eco info is being used" --------
$p17 = $d.Paragraphs(17)
$xmlG = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Summary of results</w:t></w:r><w:r><w:t>- eco info is being used</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>'
$p17.Range.InsertXML($xmlG)

$p17b = $d.Paragraphs(17)
$p17b.Range.InsertParagraphAfter()
$newBlock2 = $d.Paragraphs(18)
$xmlG2 = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Fishery&gt; physical </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>env</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &gt; ecological</w:t></w:r></w:p>' +
         '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Factors that affect observation processes are more common than ecological processes</w:t></w:r></w:p>'
$newBlock2.Range.InsertXML($xmlG2)

# --- H) Paragraph 10: "Importance: Stock status and revenue (...)" --------
$p10 = $d.Paragraphs(10)
$xmlH = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Importance: </w:t></w:r><w:r><w:t>Stock status</w:t></w:r><w:r><w:t xml:space="preserve"> (“importance”)</w:t></w:r></w:p>'
$p10.Range.InsertXML($xmlH)

Write-Output "edit complete"
